$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 247, shifting the existing
# rows 247:268 down to 249:270 (dimension grows from A1:T268 to A1:T270).
$ws.Rows("247:248").Insert()

# --- New row 247 ---
$ws.Range("A247").Value = 4
$ws.Range("B247").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C247").Value = "Los Lagos"
$ws.Range("D247").Value = 44946
$ws.Range("E247").Value = 10
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100103
$ws.Range("H247").Value = "Frutos de hueso (carozo)"
$ws.Range("I247").Value = 100103002
$ws.Range("J247").Value = "Ciruela"
$ws.Range("K247").Value = "Black Amber"
$ws.Range("L247").Value = "Primera"
$ws.Range("M247").Value = 400
$ws.Range("N247").Value = 17000
$ws.Range("O247").Value = 18000
$ws.Range("P247").Value = 17500
$ws.Range("Q247").Value = "`$/caja 14 kilos granel"
$ws.Range("R247").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("S247").Value = 1250
$ws.Range("T247").Value = 14

# --- New row 248 ---
$ws.Range("A248").Value = 4
$ws.Range("B248").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C248").Value = "Los Lagos"
$ws.Range("D248").Value = 44946
$ws.Range("E248").Value = 10
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100103
$ws.Range("H248").Value = "Frutos de hueso (carozo)"
$ws.Range("I248").Value = 100103002
$ws.Range("J248").Value = "Ciruela"
$ws.Range("K248").Value = "Lemon"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 400
$ws.Range("N248").Value = 17000
$ws.Range("O248").Value = 18000
$ws.Range("P248").Value = 17500
$ws.Range("Q248").Value = "`$/caja 14 kilos granel"
$ws.Range("R248").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("S248").Value = 1250
$ws.Range("T248").Value = 14

# Make sure the date cells keep/receive the expected date number format
# (same as the rest of column D).
$ws.Range("D247:D248").NumberFormat = "YYYY-MM-DD HH:MM:SS"
